$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 02:29"

$ws.Range("B4").Value = 7773196
$ws.Range("C4").Value = 45687
$ws.Range("D4").Value = 4975149
$ws.Range("E4").Value = 2581342
$ws.Range("G4").Value = 853
$ws.Range("H4").Value = 216705

$ws.Range("E6").Value = 460966
$ws.Range("G6").Value = 733
$ws.Range("H6").Value = 148304

$ws.Range("B11").Value = 835662
$ws.Range("C11").Value = 2733
$ws.Range("D11").Value = 723606
$ws.Range("E11").Value = 79047
$ws.Range("G11").Value = 95
$ws.Range("H11").Value = 33009

$ws.Range("B29").Value = 173123
$ws.Range("C29").Value = 1800
$ws.Range("D29").Value = 145666
$ws.Range("E29").Value = 17916
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 9541

# Row 47: Nepal -> Chequia
$ws.Range("A47").Value = "Chequia"
$ws.Range("B47").Value = 95360
$ws.Range("C47").Value = 5338
$ws.Range("D47").Value = 50767
$ws.Range("E47").Value = 43764
$ws.Range("G47").Value = 35
$ws.Range("H47").Value = 829

# Row 48: Chequia -> Nepal
$ws.Range("A48").Value = "Nepal"
$ws.Range("B48").Value = 94253
$ws.Range("C48").Value = 3439
$ws.Range("D48").Value = 68668
$ws.Range("E48").Value = 25007
$ws.Range("G48").Value = 15
$ws.Range("H48").Value = 578

$ws.Range("B69").Value = 46435
$ws.Range("C69").Value = 788
$ws.Range("D69").Value = 29270
$ws.Range("E69").Value = 16176
$ws.Range("G69").Value = 23
$ws.Range("H69").Value = 989

$ws.Range("B95").Value = 15013
$ws.Range("C95").Value = 229
$ws.Range("E95").Value = 2875

$ws.Range("E126").Value = 886
$ws.Range("H126").Value = 153

$ws.Range("B130").Value = 4979
$ws.Range("C130").Value = 14
$ws.Range("D130").Value = 4781
$ws.Range("E130").Value = 92

$ws.Range("D143").Value = 2235
$ws.Range("E143").Value = 1261

# Row 151: Guinea-Bisau -> Benin
$ws.Range("A151").Value = "Benin"
$ws.Range("B151").Value = 2411
$ws.Range("C151").Value = 54
$ws.Range("D151").Value = 1973
$ws.Range("E151").Value = 397
$ws.Range("H151").Value = 41

# Row 152: Benin -> Guinea-Bisau
$ws.Range("A152").Value = "Guinea-Bisau"
$ws.Range("B152").Value = 2385
$ws.Range("D152").Value = 1728
$ws.Range("E152").Value = 617
$ws.Range("H152").Value = 40

# Row 153: Sierra Leona -> Polinesia Francesa
$ws.Range("A153").Value = "Polinesia Francesa"
$ws.Range("B153").Value = 2358
$ws.Range("C153").Value = 130
$ws.Range("D153").Value = 1857
$ws.Range("E153").Value = 491
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 10

# Row 154: Letonia -> Sierra Leona
$ws.Range("A154").Value = "Sierra Leona"
$ws.Range("B154").Value = 2287
$ws.Range("C154").Value = 10
$ws.Range("D154").Value = 1716
$ws.Range("E154").Value = 499
$ws.Range("H154").Value = 72

# Row 155: Belice -> Letonia
$ws.Range("A155").Value = "Letonia"
$ws.Range("B155").Value = 2261
$ws.Range("C155").Value = 67
$ws.Range("D155").Value = 1322
$ws.Range("E155").Value = 899
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 40

# Row 156: Polinesia Francesa -> Belice
$ws.Range("A156").Value = "Belice"
$ws.Range("B156").Value = 2243
$ws.Range("C156").Value = 39
$ws.Range("D156").Value = 1392
$ws.Range("E156").Value = 817
$ws.Range("G156").Value = 4
$ws.Range("H156").Value = 34

$ws.Range("D167").Value = 1122
$ws.Range("E167").Value = 9

# Row 207: Santa Lucia -> Nueva Caledonia
$ws.Range("A207").Value = "Nueva Caledonia"

# Row 208: Nueva Caledonia -> Santa Lucia
$ws.Range("A208").Value = "Santa Lucia"

# Row 215: Islas Malvinas -> Montserrat
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216: Montserrat -> Islas Malvinas
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
